# Auto update Excel log
# Appends new sensor event rows to the "Proximity" and "Camera" sheets.

$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append rows 10-12 -----------------------------------
$proximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-01-28", "17:43:31", "17:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-01-28", "17:43:36", "17:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-01-28", "17:43:38", "17:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)

$startRow = 10
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $proximityRows[$i]
    # Force column A (date-looking text) to stay plain text instead of
    # being auto-converted into a date serial value by Excel.
    $proximity.Cells.Item($r, 1).NumberFormat = "@"
    $proximity.Cells.Item($r, 1).Value = $rowData[0]
    $proximity.Cells.Item($r, 2).Value = $rowData[1]
    $proximity.Cells.Item($r, 3).Value = $rowData[2]
    $proximity.Cells.Item($r, 4).Value = $rowData[3]
    $proximity.Cells.Item($r, 5).Value = $rowData[4]
    $proximity.Cells.Item($r, 6).Value = $rowData[5]
}

# --- Camera sheet: append rows 6-8 ----------------------------------------
$camera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-01-28", "17:43:32", "17:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-01-28", "17:43:35", "17:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-01-28", "17:43:39", "17:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow = 6
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $cameraRows[$i]
    $camera.Cells.Item($r, 1).NumberFormat = "@"
    $camera.Cells.Item($r, 1).Value = $rowData[0]
    $camera.Cells.Item($r, 2).Value = $rowData[1]
    $camera.Cells.Item($r, 3).Value = $rowData[2]
    $camera.Cells.Item($r, 4).Value = $rowData[3]
    $camera.Cells.Item($r, 5).Value = $rowData[4]
    $camera.Cells.Item($r, 6).Value = $rowData[5]
}
